$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 101 (pushing existing rows 101..215 down to 102..216)
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A101").Value = 3
$ws.Range("B101").Value = "Femacal de La Calera"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44781
$ws.Range("E101").Value = 5
$ws.Range("F101").Value = 100112010
$ws.Range("G101").Value = "Achicoria"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 50
$ws.Range("K101").Value = 7000
$ws.Range("L101").Value = 7000
$ws.Range("M101").Value = 7000
$ws.Range("N101").Value = "$/caja 16 unidades"
$ws.Range("O101").Value = "Provincia de Quillota"
$ws.Range("P101").Value = 438
$ws.Range("Q101").Value = 16
$ws.Range("R101").Value = "Hortaliza"
